$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant columns (same across all data rows 2-30): A,B,C,E,F,G,H,I,J,K
$mercadoId = 4
$mercado = "Feria Lagunitas de Puerto Montt"
$region = "Los Lagos"
$codreg = 10
$tipo = "Fruta"
$productoId = 100107
$producto = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad = "Cultivar IV Región"

# Row-specific data: D(fecha serial), L(calidad), M(volumen), N(min), O(max), P(prom), Q(unidad), R(origen), S(precio/kg), T(kg/unidad)
$rows = @(
    @(44530, "Primera", 200, 19000, 20000, 19500, "`$/bandeja 8 kilos", "Provincia de Limarí", 2438, 8),
    @(44530, "Segunda", 100, 16000, 16000, 16000, "`$/bandeja 8 kilos", "Provincia de Limarí", 2000, 8),
    @(44512, "Segunda", 300, 19000, 20000, 19500, "`$/bandeja 8 kilos", "Provincia de Limarí", 2438, 8),
    @(45247, "Primera", 100, 25000, 25000, 25000, "`$/bandeja 10 kilos", "Provincia del Elquí", 2500, 10),
    @(44533, "Primera", 300, 18000, 19000, 18500, "`$/bandeja 8 kilos", "Provincia de Limarí", 2312, 8),
    @(44533, "Segunda", 100, 16000, 16000, 16000, "`$/bandeja 8 kilos", "Provincia de Limarí", 2050, 8),
    @(44523, "Primera", 400, 21000, 22000, 21500, "`$/bandeja 8 kilos", "Provincia de Limarí", 2688, 8),
    @(44523, "Segunda", 100, 18000, 18000, 18000, "`$/bandeja 8 kilos", "Provincia de Limarí", 2250, 8),
    @(44162, "Primera", 200, 2000, 2100, 2050, "`$/kilo (en caja de 14 kilos)", "Provincia de Limarí", 2050, 1),
    @(45240, "Primera", 150, 28000, 28000, 28000, "`$/bandeja 10 kilos", "Provincia del Elquí", 2800, 10),
    @(44519, "Primera", 400, 21000, 22500, 21500, "`$/bandeja 8 kilos", "Provincia de Limarí", 2688, 8),
    @(44519, "Segunda", 200, 18000, 18000, 18000, "`$/bandeja 8 kilos", "Provincia de Limarí", 2250, 8),
    @(44498, "Segunda", 300, 19000, 20000, 19500, "`$/bandeja 8 kilos", "Provincia de Limarí", 2438, 8),
    @(44516, "Segunda", 200, 18000, 19000, 18500, "`$/bandeja 8 kilos", "Provincia de Limarí", 2312, 8),
    @(44890, "Primera", 200, 22000, 22500, 22250, "`$/bandeja 8 kilos", "Provincia de Limarí", 2781, 8),
    @(44495, "Segunda", 270, 19000, 20000, 19556, "`$/bandeja 8 kilos", "Provincia de Limarí", 2444, 8),
    @(44873, "Primera", 300, 22000, 22500, 22250, "`$/bandeja 8 kilos", "Provincia de Limarí", 2781, 8),
    @(44505, "Segunda", 300, 19000, 20000, 19500, "`$/bandeja 8 kilos", "Provincia de Limarí", 2438, 8),
    @(44159, "Primera", 300, 2000, 2100, 2050, "`$/kilo (en caja de 14 kilos)", "Provincia de Limarí", 2050, 1),
    @(45237, "Primera", 150, 28000, 28000, 28000, "`$/bandeja 10 kilos", "Provincia del Elquí", 2800, 10),
    @(44526, "Primera", 300, 21000, 21000, 21000, "`$/bandeja 8 kilos", "Provincia de Limarí", 2625, 8),
    @(44488, "Segunda", 160, 17000, 18000, 17500, "`$/bandeja 8 kilos", "Provincia de Limarí", 2188, 8),
    @(45244, "Primera", 150, 25000, 25000, 25000, "`$/bandeja 10 kilos", "Provincia del Elquí", 2500, 10),
    @(44876, "Primera", 300, 22000, 22500, 22250, "`$/bandeja 8 kilos", "Provincia de Limarí", 2781, 8),
    @(44491, "Segunda", 200, 18000, 19000, 18500, "`$/bandeja 8 kilos", "Provincia de Limarí", 2312, 8),
    @(44880, "Primera", 300, 22000, 22500, 22250, "`$/bandeja 8 kilos", "Provincia de Limarí", 2781, 8),
    @(44509, "Segunda", 200, 19000, 20000, 19500, "`$/bandeja 8 kilos", "Provincia de Limarí", 2438, 8),
    @(44894, "Primera", 200, 22000, 22500, 22250, "`$/bandeja 8 kilos", "Provincia de Limarí", 2781, 8),
    @(44895, "Primera", 200, 22000, 22500, 22250, "`$/bandeja 8 kilos", "Provincia de Limarí", 2781, 8),
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $row[0]
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad

    $ws.Cells.Item($r, 12).Value = $row[1]
    $ws.Cells.Item($r, 13).Value = $row[2]
    $ws.Cells.Item($r, 14).Value = $row[3]
    $ws.Cells.Item($r, 15).Value = $row[4]
    $ws.Cells.Item($r, 16).Value = $row[5]
    $ws.Cells.Item($r, 17).Value = $row[6]
    $ws.Cells.Item($r, 18).Value = $row[7]
    $ws.Cells.Item($r, 19).Value = $row[8]
    $ws.Cells.Item($r, 20).Value = $row[9]
}

Write-Host "Final dimension: $($ws.UsedRange.Address())"
